# swift.docx edit: rewrite the opening sentence's run and append a
# red "(This is a change - Version for main branch)" note, split
# across three runs (as captured by the target OOXML diff).

$d = $word.ActiveDocument

# --- paragraph 1: "This is a Microsoft word document." -------------------
$para1 = $d.Paragraphs(1).Range

# Work only inside the paragraph (exclude the trailing paragraph mark)
# so the paragraph's own identity/rsid attributes are left untouched.
$body = $d.Range($para1.Start, $para1.End - 1)

# Rewrite the existing single run's text, adding two trailing spaces
# (Word will emit xml:space="preserve" automatically because of them).
$body.Text = "This is a Microsoft word document.  "

# Move a collapsed insertion point to the end of what we just wrote.
$body.Collapse(0)

# Run 2 (red): "(This is a change " + EN DASH + " Ve"
$run2 = $body.Duplicate
$run2.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$run2.Font.Color = 255

# Run 3 (red): "rsion for main branch"
$run3 = $run2.Duplicate
$run3.Collapse(0)
$run3.InsertAfter("rsion for main branch")
$run3.Font.Color = 255

# Run 4 (red): ")"
$run4 = $run3.Duplicate
$run4.Collapse(0)
$run4.InsertAfter(")")
$run4.Font.Color = 255
